# ---------------------------------------------------------------------------
# Add a new "New Requirement" worksheet (with the new-column requirements
# table) at the end of the workbook, make it the active sheet, and tidy up
# the selections left on the other two sheets.
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook

$wsData       = $wb.Worksheets.Item(1)
$wsValidation = $wb.Worksheets.Item(2)

# --- create the new sheet as the last tab ----------------------------------
$wsNew = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsNew.Name = "New Requirement"

# --- instructional notes (rows 3-5) -----------------------------------------
$wsNew.Range("A3").Value = "Edit the file data/file template.xls to add the below new columns based on table below"
$wsNew.Range("A4").Value = "use File Mandatory column from below table to change column background column to grey for columns defined as mandatory"
$wsNew.Range("A5").Value = "use Sequence column from below table to decide the insertion point of column"

# --- header row (row 7) ------------------------------------------------------
$wsNew.Range("A7").Value = "Attribute Name"
$wsNew.Range("B7").Value = "Sequence "
$wsNew.Range("C7").Value = "File Mandatory"
$wsNew.Range("D7").Value = "Target Mandatory"
$wsNew.Range("E7").Value = "Field Description"
$wsNew.Range("F7").Value = "Sample Values"
$wsNew.Range("G7").Value = "Comments"

# --- data rows (8-10) --------------------------------------------------------
$wsNew.Range("A8").Value = "Adjustment Type"
$wsNew.Range("B8").Value = 3
$wsNew.Range("C8").Value = "O - Optional"
$wsNew.Range("D8").Value = "O - Optional"
$wsNew.Range("E8").Value = "explains the adjustment type based on adjusment definition"
$wsNew.Range("F8").Value = "Adding/Removing  Records, Overwrite Missing, Manual Data Transfer"
$wsNew.Range("G8").Value = "a drop down list will be available for user while filling"

$wsNew.Range("A9").Value = "Frequency"
$wsNew.Range("B9").Value = 1
$wsNew.Range("C9").Value = "M - Mandatory"
$wsNew.Range("D9").Value = "M - Mandatory"
$wsNew.Range("E9").Value = "frequency of adjustment"
$wsNew.Range("F9").Value = "Daily, Monthly, Quarterly, Yearly, Ad-hoc"
$wsNew.Range("G9").Value = "a drop down list will be available for user while filling"

$wsNew.Range("A10").Value = "Reason CodeLevel 1"
$wsNew.Range("B10").Value = 5
$wsNew.Range("C10").Value = "CM - Conditional Mandatory"
$wsNew.Range("D10").Value = "M - Mandatory"
$wsNew.Range("E10").Value = "Leve 1 reason for adjustment"

# --- formatting ---------------------------------------------------------------
# Header row: copy the orange "bold header" look used on the Validation sheet,
# then centre it (this produces the extra centred cellXfs entry).
$wsValidation.Range("A1:D1").Copy() | Out-Null
$wsNew.Range("A7:G7").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$wsNew.Range("A7:G7").HorizontalAlignment = -4108      # xlCenter
$wsNew.Application.CutCopyMode = $false

# Wrap the long "Adding/Removing..." cell and grow its row to fit two lines.
$wsNew.Range("F8").WrapText = $true
$wsNew.Rows.Item(8).RowHeight = 28.8

# Column widths approximating the authored sheet.
$wsNew.Columns.Item(1).ColumnWidth = 14.6640625
$wsNew.Columns.Item(2).ColumnWidth = 14.6640625
$wsNew.Columns.Item(3).ColumnWidth = 23.21875
$wsNew.Columns.Item(4).ColumnWidth = 14.6640625
$wsNew.Columns.Item(5).ColumnWidth = 23.77734375
$wsNew.Columns.Item(6).ColumnWidth = 33.33203125
$wsNew.Columns.Item(7).ColumnWidth = 42.5546875

# --- selections / active sheet -------------------------------------------------
$wsValidation.Activate()
$wsValidation.Range("B1:C1").Select()

$wsNew.Activate()
$wsNew.Range("B10").Select()
